# Auto-generated Excel COM-interop script applying the scraped-data refresh
# for sheet "LP1912" (rows 2,3 header text + row reorder/refresh 50-500),
# plus the "Ultima actualizacion" timestamp bump on the other two sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LP1912")

$ws.Cells.Item(2, 1).Value = "Última actualización: 20:52:54"
$ws.Cells.Item(3, 1).Value = "Total filas: 495"
$ws.Cells.Item(50, 1).Value = "07:12:53"
$ws.Cells.Item(50, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(50, 4).Value = 49
$ws.Cells.Item(51, 1).Value = "06:45:50"
$ws.Cells.Item(51, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(51, 4).Value = 76
$ws.Cells.Item(64, 1).Value = "08:29:19"
$ws.Cells.Item(64, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(65, 1).Value = "06:45:50"
$ws.Cells.Item(65, 3).Value = "14_ABASTO"
$ws.Cells.Item(65, 4).Value = 104
$ws.Cells.Item(81, 1).Value = "07:36:59"
$ws.Cells.Item(81, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(81, 4).Value = 86
$ws.Cells.Item(82, 1).Value = "07:12:53"
$ws.Cells.Item(82, 3).Value = "17X38_ROMERO"
$ws.Cells.Item(82, 4).Value = 110
$ws.Cells.Item(206, 1).Value = "10:55:25"
$ws.Cells.Item(206, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(206, 4).Value = 108
$ws.Cells.Item(207, 1).Value = "10:48:14"
$ws.Cells.Item(207, 3).Value = "14_ABASTO"
$ws.Cells.Item(207, 4).Value = 115
$ws.Cells.Item(220, 1).Value = "12:11:45"
$ws.Cells.Item(220, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(220, 4).Value = 55
$ws.Cells.Item(221, 1).Value = "11:11:31"
$ws.Cells.Item(221, 3).Value = "14_ABASTO"
$ws.Cells.Item(221, 4).Value = 115
$ws.Cells.Item(222, 1).Value = "11:53:59"
$ws.Cells.Item(222, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(222, 4).Value = 74
$ws.Cells.Item(223, 1).Value = "11:47:13"
$ws.Cells.Item(223, 3).Value = "14_ABASTO"
$ws.Cells.Item(223, 4).Value = 80
$ws.Cells.Item(238, 1).Value = "11:34:25"
$ws.Cells.Item(238, 3).Value = "215A_EL PATO"
$ws.Cells.Item(238, 4).Value = 116
$ws.Cells.Item(239, 1).Value = "11:53:59"
$ws.Cells.Item(239, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(239, 4).Value = 97
$ws.Cells.Item(240, 1).Value = "11:47:13"
$ws.Cells.Item(240, 3).Value = "10_OLMOS"
$ws.Cells.Item(240, 4).Value = 103
$ws.Cells.Item(283, 1).Value = "12:45:57"
$ws.Cells.Item(283, 3).Value = "215B_EL PATO"
$ws.Cells.Item(283, 4).Value = 119
$ws.Cells.Item(284, 1).Value = "14:44:53"
$ws.Cells.Item(284, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(284, 4).Value = 0
$ws.Cells.Item(313, 1).Value = "14:10:21"
$ws.Cells.Item(313, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(313, 4).Value = 107
$ws.Cells.Item(314, 1).Value = "14:44:53"
$ws.Cells.Item(314, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(314, 4).Value = 73
$ws.Cells.Item(376, 1).Value = "16:13:19"
$ws.Cells.Item(376, 3).Value = "17_ROMERO"
$ws.Cells.Item(376, 4).Value = 110
$ws.Cells.Item(377, 1).Value = "17:34:55"
$ws.Cells.Item(377, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(377, 4).Value = 29
$ws.Cells.Item(406, 3).Value = "15_ABASTO"
$ws.Cells.Item(407, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(422, 1).Value = "18:44:14"
$ws.Cells.Item(422, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(422, 4).Value = 33
$ws.Cells.Item(423, 1).Value = "18:10:23"
$ws.Cells.Item(423, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(423, 4).Value = 67
$ws.Cells.Item(480, 1).Value = "20:52:54"
$ws.Cells.Item(480, 2).Value = "20:52"
$ws.Cells.Item(480, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(480, 4).Value = 0
$ws.Cells.Item(481, 1).Value = "19:47:42"
$ws.Cells.Item(481, 2).Value = "20:55"
$ws.Cells.Item(481, 4).Value = 68
$ws.Cells.Item(482, 2).Value = "20:56"
$ws.Cells.Item(482, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(482, 4).Value = 105
$ws.Cells.Item(483, 1).Value = "20:52:54"
$ws.Cells.Item(483, 2).Value = "20:56"
$ws.Cells.Item(483, 3).Value = "17X38_ROMERO"
$ws.Cells.Item(483, 4).Value = 4
$ws.Cells.Item(484, 1).Value = "19:11:56"
$ws.Cells.Item(484, 2).Value = "21:01"
$ws.Cells.Item(484, 3).Value = "215A_EL PATO"
$ws.Cells.Item(484, 4).Value = 110
$ws.Cells.Item(485, 1).Value = "19:35:19"
$ws.Cells.Item(485, 2).Value = "21:02"
$ws.Cells.Item(485, 3).Value = "215A_EL PATO"
$ws.Cells.Item(485, 4).Value = 87
$ws.Cells.Item(486, 1).Value = "20:45:44"
$ws.Cells.Item(486, 2).Value = "21:06"
$ws.Cells.Item(486, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(486, 4).Value = 21
$ws.Cells.Item(487, 1).Value = "19:47:42"
$ws.Cells.Item(487, 2).Value = "21:09"
$ws.Cells.Item(487, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(487, 4).Value = 82
$ws.Cells.Item(488, 1).Value = "19:47:42"
$ws.Cells.Item(488, 2).Value = "21:23"
$ws.Cells.Item(488, 3).Value = "10_OLMOS"
$ws.Cells.Item(488, 4).Value = 96
$ws.Cells.Item(489, 1).Value = "19:35:19"
$ws.Cells.Item(489, 2).Value = "21:24"
$ws.Cells.Item(489, 3).Value = "10_OLMOS"
$ws.Cells.Item(489, 4).Value = 109
$ws.Cells.Item(490, 1).Value = "20:52:54"
$ws.Cells.Item(490, 2).Value = "21:29"
$ws.Cells.Item(490, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(490, 4).Value = 37
$ws.Cells.Item(491, 1).Value = "20:45:44"
$ws.Cells.Item(491, 2).Value = "21:30"
$ws.Cells.Item(491, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(491, 4).Value = 45
$ws.Cells.Item(492, 1).Value = "19:54:54"
$ws.Cells.Item(492, 2).Value = "21:48"
$ws.Cells.Item(492, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(492, 4).Value = 114
$ws.Cells.Item(493, 1).Value = "20:11:44"
$ws.Cells.Item(493, 2).Value = "21:49"
$ws.Cells.Item(493, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(493, 4).Value = 98
$ws.Cells.Item(494, 1).Value = "20:11:44"
$ws.Cells.Item(494, 2).Value = "21:55"
$ws.Cells.Item(494, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(494, 4).Value = 104
$ws.Cells.Item(495, 2).Value = "22:18"
$ws.Cells.Item(495, 3).Value = "10_OLMOS"
$ws.Cells.Item(495, 4).Value = 107
$ws.Cells.Item(496, 1).Value = "20:31:53"
$ws.Cells.Item(496, 2).Value = "22:25"
$ws.Cells.Item(496, 3).Value = "15_ABASTO"
$ws.Cells.Item(496, 4).Value = 114
$ws.Cells.Item(496, 5).Value = "LP1912"
$ws.Cells.Item(497, 1).Value = "20:31:53"
$ws.Cells.Item(497, 2).Value = "22:29"
$ws.Cells.Item(497, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(497, 4).Value = 118
$ws.Cells.Item(497, 5).Value = "LP1912"
$ws.Cells.Item(498, 1).Value = "20:31:53"
$ws.Cells.Item(498, 2).Value = "22:30"
$ws.Cells.Item(498, 3).Value = "215C_EL PATO"
$ws.Cells.Item(498, 4).Value = 119
$ws.Cells.Item(498, 5).Value = "LP1912"
$ws.Cells.Item(499, 1).Value = "20:52:54"
$ws.Cells.Item(499, 2).Value = "22:35"
$ws.Cells.Item(499, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(499, 4).Value = 103
$ws.Cells.Item(499, 5).Value = "LP1912"
$ws.Cells.Item(500, 1).Value = "20:52:54"
$ws.Cells.Item(500, 2).Value = "22:48"
$ws.Cells.Item(500, 3).Value = "14_ABASTO"
$ws.Cells.Item(500, 4).Value = 116
$ws.Cells.Item(500, 5).Value = "LP1912"

# Other two sheets only get their "Ultima actualizacion" timestamp bumped.
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2, 1).Value = "Última actualización: 20:52:54"

$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2, 1).Value = "Última actualización: 20:52:54"

